$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-11-08 Saturday" "2025-11-09 Sunday"

Replace-Text "847÷7=121, 0" "886÷3=295, 1"
Replace-Text "856÷5=171, 1" "669÷6=111, 3"
Replace-Text "683÷5=136, 3" "635÷8=79, 3"
Replace-Text "664÷7=94, 6" "978÷8=122, 2"
Replace-Text "175÷5=35, 0" "123÷9=13, 6"

Replace-Text "795÷2=397, 1" "121÷3=40, 1"
Replace-Text "591÷7=84, 3" "524÷3=174, 2"
Replace-Text "660÷9=73, 3" "174÷3=58, 0"
Replace-Text "971÷9=107, 8" "586÷6=97, 4"
Replace-Text "656÷4=164, 0" "518÷6=86, 2"

Replace-Text "943÷6=157, 1" "684÷9=76, 0"
Replace-Text "530÷2=265, 0" "691÷5=138, 1"
Replace-Text "132÷5=26, 2" "172÷6=28, 4"
Replace-Text "253÷6=42, 1" "797÷6=132, 5"
Replace-Text "948÷7=135, 3" "301÷6=50, 1"

Replace-Text "762÷2=381, 0" "374÷2=187, 0"
Replace-Text "456÷4=114, 0" "365÷4=91, 1"
Replace-Text "497÷5=99, 2" "506÷5=101, 1"
Replace-Text "797÷5=159, 2" "387÷5=77, 2"
Replace-Text "201÷3=67, 0" "390÷8=48, 6"

Replace-Text "540÷5=108, 0" "610÷7=87, 1"
Replace-Text "666÷8=83, 2" "933÷2=466, 1"
Replace-Text "795÷9=88, 3" "211÷8=26, 3"
Replace-Text "256÷4=64, 0" "569÷4=142, 1"
Replace-Text "510÷6=85, 0" "770÷3=256, 2"

Write-Output "Done"
